$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row138
$ws.Range("H138").Value = 2213.138
$ws.Range("I138").Value = 2173
$ws.Range("J138").Value = 2228.4285
$ws.Range("K138").Value = 6519
$ws.Range("L138").Value = 6685.2855
$ws.Range("M138").Value = -1379
$ws.Range("N138").Value = -16965.2855

$ws = $wb.Worksheets.Item("ARM")
# ARM!row2
$ws.Range("H2").Value = 1759.909
$ws.Range("I2").Value = 1531.4286
$ws.Range("J2").Value = 2159.75
$ws.Range("K2").Value = 1531.4286
$ws.Range("L2").Value = 2159.75
$ws.Range("M2").Value = -1418.4286
$ws.Range("N2").Value = -2385.75

# ARM!row74
$ws.Range("H74").Value = 2854.5715
$ws.Range("I74").Value = 4190
$ws.Range("J74").Value = 2112.6667
$ws.Range("K74").Value = 4190
$ws.Range("L74").Value = 2112.6667
$ws.Range("M74").Value = -3316
$ws.Range("N74").Value = -3860.6667

# ARM!row77
$ws.Range("H77").Value = 2854.5715
$ws.Range("I77").Value = 4190
$ws.Range("J77").Value = 2112.6667
$ws.Range("K77").Value = 20950
$ws.Range("L77").Value = 10563.3335
$ws.Range("M77").Value = -16582
$ws.Range("N77").Value = -19299.3335

# ARM!row116
$ws.Range("H116").Value = 1759.909
$ws.Range("I116").Value = 1531.4286
$ws.Range("J116").Value = 2159.75
$ws.Range("K116").Value = 1531.4286
$ws.Range("L116").Value = 2159.75
$ws.Range("M116").Value = 762.5714
$ws.Range("N116").Value = -6747.75

# ARM!row117
$ws.Range("H117").Value = 19248
$ws.Range("J117").Value = 19248
$ws.Range("L117").Value = 19248
$ws.Range("N117").Value = -28426

# ARM!row132
$ws.Range("H132").Value = 2558.8484
$ws.Range("I132").Value = 1173.7894
$ws.Range("J132").Value = 4438.5713
$ws.Range("K132").Value = 3521.3682
$ws.Range("L132").Value = 13315.7139
$ws.Range("M132").Value = -991.3681999999999
$ws.Range("N132").Value = -18375.7139

$ws = $wb.Worksheets.Item("BSM")
# BSM!row3
$ws.Range("H3").Value = 1759.909
$ws.Range("I3").Value = 1531.4286
$ws.Range("J3").Value = 2159.75
$ws.Range("K3").Value = 1531.4286
$ws.Range("L3").Value = 2159.75
$ws.Range("M3").Value = -1417.4286
$ws.Range("N3").Value = -2387.75

# BSM!row20
$ws.Range("H20").Value = 2813.7368
$ws.Range("I20").Value = 3458
$ws.Range("J20").Value = 1927.875
$ws.Range("K20").Value = 3458
$ws.Range("L20").Value = 1927.875
$ws.Range("M20").Value = -3211
$ws.Range("N20").Value = -2421.875

# BSM!row139
$ws.Range("H139").Value = 47980
$ws.Range("J139").Value = 47980
$ws.Range("L139").Value = 47980
$ws.Range("N139").Value = -58260

$ws = $wb.Worksheets.Item("CRP")
# CRP!row31
$ws.Range("H31").Value = 21858.846
$ws.Range("I31").Value = 2206
$ws.Range("J31").Value = 25432.092
$ws.Range("K31").Value = 2206
$ws.Range("L31").Value = 25432.092
$ws.Range("M31").Value = -1911
$ws.Range("N31").Value = -26022.092

# CRP!row34
$ws.Range("H34").Value = 21858.846
$ws.Range("I34").Value = 2206
$ws.Range("J34").Value = 25432.092
$ws.Range("K34").Value = 2206
$ws.Range("L34").Value = 25432.092
$ws.Range("M34").Value = -2004
$ws.Range("N34").Value = -25836.092

# CRP!row134
$ws.Range("H134").Value = 6862.0713
$ws.Range("I134").Value = 6171.25
$ws.Range("K134").Value = 18513.75
$ws.Range("M134").Value = -15978.75

$ws = $wb.Worksheets.Item("CUL")
# CUL!row5
$ws.Range("H5").Value = 1089.8846
$ws.Range("I5").Value = 526.7059
$ws.Range("J5").Value = 2153.6667
$ws.Range("K5").Value = 1580.1177
$ws.Range("L5").Value = 6461.000100000001
$ws.Range("M5").Value = -1468.1177
$ws.Range("N5").Value = -6685.000100000001

# CUL!row135
$ws.Range("H135").Value = 1089.8846
$ws.Range("I135").Value = 526.7059
$ws.Range("J135").Value = 2153.6667
$ws.Range("K135").Value = 4740.3531
$ws.Range("L135").Value = 19383.0003
$ws.Range("M135").Value = -2205.3531
$ws.Range("N135").Value = -24453.0003

$ws = $wb.Worksheets.Item("GSM")
# GSM!row12
$ws.Range("H12").Value = 6125625.5
$ws.Range("I12").Value = 6000714.5
$ws.Range("J12").Value = 7000000
$ws.Range("K12").Value = 6000714.5
$ws.Range("L12").Value = 7000000
$ws.Range("M12").Value = -6000574.5
$ws.Range("N12").Value = -7000280

# GSM!row18
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# GSM!row43
$ws.Range("H43").Value = 2000
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

# GSM!row46
$ws.Range("H46").Value = 7810
$ws.Range("I46").Value = 3500
$ws.Range("J46").Value = 9965
$ws.Range("K46").Value = 3500
$ws.Range("L46").Value = 9965
$ws.Range("M46").Value = -3344
$ws.Range("N46").Value = -10277

# GSM!row70
$ws.Range("H70").Value = 9046
$ws.Range("I70").Value = 10728
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 10728
$ws.Range("L70").Value = 4000
$ws.Range("M70").Value = -10458
$ws.Range("N70").Value = -4540

# GSM!row73
$ws.Range("H73").Value = 9046
$ws.Range("I73").Value = 10728
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 10728
$ws.Range("L73").Value = 4000
$ws.Range("M73").Value = -9792
$ws.Range("N73").Value = -5872

# GSM!row80
$ws.Range("H80").Value = 3177
$ws.Range("I80").Value = 3333
$ws.Range("J80").Value = 3151
$ws.Range("K80").Value = 3333
$ws.Range("L80").Value = 3151
$ws.Range("M80").Value = -2335
$ws.Range("N80").Value = -5147

# GSM!row83
$ws.Range("H83").Value = 3177
$ws.Range("I83").Value = 3333
$ws.Range("J83").Value = 3151
$ws.Range("K83").Value = 16665
$ws.Range("L83").Value = 15755
$ws.Range("M83").Value = -11673
$ws.Range("N83").Value = -25739

# GSM!row113
$ws.Range("H113").Value = 1907.3334
$ws.Range("I113").Value = 1519.125
$ws.Range("J113").Value = 5013
$ws.Range("K113").Value = 1519.125
$ws.Range("L113").Value = 5013
$ws.Range("M113").Value = 650.875
$ws.Range("N113").Value = -9353

# GSM!row126
$ws.Range("H126").Value = 2109.4614
$ws.Range("I126").Value = 1678.4375
$ws.Range("J126").Value = 2799.1
$ws.Range("K126").Value = 5035.3125
$ws.Range("L126").Value = 8397.299999999999
$ws.Range("M126").Value = -2565.3125
$ws.Range("N126").Value = -13337.3

$ws = $wb.Worksheets.Item("LTW")
# LTW!row61
$ws.Range("H61").Value = 3100
$ws.Range("I61").Value = 2200
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 2200
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -1998
$ws.Range("N61").Value = -4404

# LTW!row113
$ws.Range("H113").Value = 3100
$ws.Range("I113").Value = 2200
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 2200
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -30
$ws.Range("N113").Value = -8340

# LTW!row136
$ws.Range("H136").Value = 3947.077
$ws.Range("I136").Value = 3232.842
$ws.Range("J136").Value = 5885.7144
$ws.Range("K136").Value = 9698.526
$ws.Range("L136").Value = 17657.1432
$ws.Range("M136").Value = -7148.526
$ws.Range("N136").Value = -22757.1432

# LTW!row139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
